$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 78.01006965174119
$ws.Cells.Item(2, 2).Value = 65.94399999999997
$ws.Cells.Item(2, 3).Value = 87.53200000000002
$ws.Cells.Item(3, 1).Value = 26.50179104477606
$ws.Cells.Item(3, 2).Value = 15.29600000000001
$ws.Cells.Item(3, 3).Value = 40.00399999999999
$ws.Cells.Item(4, 1).Value = 40.22799999999992
$ws.Cells.Item(4, 2).Value = 25.26400000000001
$ws.Cells.Item(4, 3).Value = 56.628
$ws.Cells.Item(5, 1).Value = 30.03992039800988
$ws.Cells.Item(5, 2).Value = 17.21199999999999
$ws.Cells.Item(5, 3).Value = 46.09600000000002
$ws.Cells.Item(6, 1).Value = 61.96770149253724
$ws.Cells.Item(6, 2).Value = 45.01999999999999
$ws.Cells.Item(6, 3).Value = 75.92000000000004
$ws.Cells.Item(7, 1).Value = 74.98384079601981
$ws.Cells.Item(7, 2).Value = 62.54399999999999
$ws.Cells.Item(7, 3).Value = 85.69999999999995
$ws.Cells.Item(8, 1).Value = 69.95639800995012
$ws.Cells.Item(8, 2).Value = 54.03199999999999
$ws.Cells.Item(8, 3).Value = 83.21199999999999
$ws.Cells.Item(9, 1).Value = 59.40698507462682
$ws.Cells.Item(9, 2).Value = 42.876
$ws.Cells.Item(9, 3).Value = 74.98800000000004
$ws.Cells.Item(10, 1).Value = 69.16362189054713
$ws.Cells.Item(10, 2).Value = 54.46
$ws.Cells.Item(10, 3).Value = 81.78799999999997
$ws.Cells.Item(11, 1).Value = 78.60533333333328
$ws.Cells.Item(11, 2).Value = 67.07199999999997
$ws.Cells.Item(11, 3).Value = 87.81599999999997
$ws.Cells.Item(12, 1).Value = 73.90105472636797
$ws.Cells.Item(12, 2).Value = 60.91600000000005
$ws.Cells.Item(12, 3).Value = 83.708
$ws.Cells.Item(13, 1).Value = 74.11414925373118
$ws.Cells.Item(13, 2).Value = 61.51599999999996
$ws.Cells.Item(13, 3).Value = 86.42799999999995
$ws.Cells.Item(14, 1).Value = 77.40965174129343
$ws.Cells.Item(14, 2).Value = 65.59999999999999
$ws.Cells.Item(14, 3).Value = 87.47999999999998
$ws.Cells.Item(15, 1).Value = 33.90991044776108
$ws.Cells.Item(15, 2).Value = 20.41599999999998
$ws.Cells.Item(15, 3).Value = 50.29200000000002
$ws.Cells.Item(16, 1).Value = 61.83279601990043
$ws.Cells.Item(16, 2).Value = 44.79999999999998
$ws.Cells.Item(16, 3).Value = 75.96400000000003
$ws.Cells.Item(17, 1).Value = 73.6232835820894
$ws.Cells.Item(17, 2).Value = 58.53600000000005
$ws.Cells.Item(17, 3).Value = 84.25599999999994
$ws.Cells.Item(18, 1).Value = 40.19564179104471
$ws.Cells.Item(18, 2).Value = 25.38400000000001
$ws.Cells.Item(18, 3).Value = 56.48400000000001
$ws.Cells.Item(19, 1).Value = 78.68248756218901
$ws.Cells.Item(19, 2).Value = 66.96799999999996
$ws.Cells.Item(19, 3).Value = 87.69199999999995
$ws.Cells.Item(20, 1).Value = 64.82710447761187
$ws.Cells.Item(20, 2).Value = 48.09599999999998
$ws.Cells.Item(20, 3).Value = 77.47200000000007
$ws.Cells.Item(21, 1).Value = 48.66288557213925
$ws.Cells.Item(21, 2).Value = 30.79999999999999
$ws.Cells.Item(21, 3).Value = 64.83599999999997
$ws.Cells.Item(22, 1).Value = 69.66310447761188
$ws.Cells.Item(22, 2).Value = 53.636
$ws.Cells.Item(22, 3).Value = 82.26400000000004
$ws.Cells.Item(23, 1).Value = 36.49297512437792
$ws.Cells.Item(23, 2).Value = 21.82000000000001
$ws.Cells.Item(23, 3).Value = 52.11200000000001
$ws.Cells.Item(24, 1).Value = 61.77723383084571
$ws.Cells.Item(24, 2).Value = 44.59599999999998
$ws.Cells.Item(24, 3).Value = 75.98400000000004
$ws.Cells.Item(25, 1).Value = 56.27596019900488
$ws.Cells.Item(25, 2).Value = 39.39599999999998
$ws.Cells.Item(25, 3).Value = 72.26799999999994
$ws.Cells.Item(26, 1).Value = 48.98340298507456
$ws.Cells.Item(26, 2).Value = 33.09599999999998
$ws.Cells.Item(26, 3).Value = 65.56000000000004
$ws.Cells.Item(27, 1).Value = 74.05108457711427
$ws.Cells.Item(27, 2).Value = 62.44000000000003
$ws.Cells.Item(27, 3).Value = 86.07599999999994
$ws.Cells.Item(28, 1).Value = 76.14407960198993
$ws.Cells.Item(28, 2).Value = 62.508
$ws.Cells.Item(28, 3).Value = 87.34400000000005
$ws.Cells.Item(29, 1).Value = 37.90326368159195
$ws.Cells.Item(29, 2).Value = 23.432
$ws.Cells.Item(29, 3).Value = 53.98399999999999
$ws.Cells.Item(30, 1).Value = 72.54756218905459
$ws.Cells.Item(30, 2).Value = 57.504
$ws.Cells.Item(30, 3).Value = 84.248
$ws.Cells.Item(31, 1).Value = 76.99124378109434
$ws.Cells.Item(31, 2).Value = 64.78000000000003
$ws.Cells.Item(31, 3).Value = 86.65999999999997
$ws.Cells.Item(32, 1).Value = 73.64163184079584
$ws.Cells.Item(32, 2).Value = 60.61199999999997
$ws.Cells.Item(32, 3).Value = 85.34800000000001
$ws.Cells.Item(33, 1).Value = 77.49112437810936
$ws.Cells.Item(33, 2).Value = 66.036
$ws.Cells.Item(33, 3).Value = 87.76799999999994
$ws.Cells.Item(34, 1).Value = 79.18712437810936
$ws.Cells.Item(34, 2).Value = 69.90000000000001
$ws.Cells.Item(34, 3).Value = 88.11999999999996
$ws.Cells.Item(35, 1).Value = 63.05259701492534
$ws.Cells.Item(35, 2).Value = 45.632
$ws.Cells.Item(35, 3).Value = 77.068
$ws.Cells.Item(36, 1).Value = 74.85148258706461
$ws.Cells.Item(36, 2).Value = 61.23600000000003
$ws.Cells.Item(36, 3).Value = 86.25199999999995
$ws.Cells.Item(37, 1).Value = 51.86445771144265
$ws.Cells.Item(37, 2).Value = 34.032
$ws.Cells.Item(37, 3).Value = 67.07999999999997
$ws.Cells.Item(38, 1).Value = 78.02796019900489
$ws.Cells.Item(38, 2).Value = 65.54000000000001
$ws.Cells.Item(38, 3).Value = 87.69199999999999
$ws.Cells.Item(39, 1).Value = 75.48875621890544
$ws.Cells.Item(39, 2).Value = 64.50799999999998
$ws.Cells.Item(39, 3).Value = 86.23999999999995
$ws.Cells.Item(40, 1).Value = 45.01876616915414
$ws.Cells.Item(40, 2).Value = 28.89600000000001
$ws.Cells.Item(40, 3).Value = 63.54000000000005
$ws.Cells.Item(41, 1).Value = 73.06989054726353
$ws.Cells.Item(41, 2).Value = 56.992
$ws.Cells.Item(41, 3).Value = 84.85199999999993
$ws.Cells.Item(42, 1).Value = 78.8741293532338
$ws.Cells.Item(42, 2).Value = 66.70399999999999
$ws.Cells.Item(42, 3).Value = 87.94399999999995
$ws.Cells.Item(43, 1).Value = 72.51998009950236
$ws.Cells.Item(43, 2).Value = 57.636
$ws.Cells.Item(43, 3).Value = 84.18400000000001
$ws.Cells.Item(44, 1).Value = 76.19601990049739
$ws.Cells.Item(44, 2).Value = 62.18000000000004
$ws.Cells.Item(44, 3).Value = 86.34
$ws.Cells.Item(45, 1).Value = 78.80881592039795
$ws.Cells.Item(45, 2).Value = 67.48800000000003
$ws.Cells.Item(45, 3).Value = 87.95599999999996
$ws.Cells.Item(46, 1).Value = 76.96577114427845
$ws.Cells.Item(46, 2).Value = 64.44400000000005
$ws.Cells.Item(46, 3).Value = 87.40399999999995
$ws.Cells.Item(47, 1).Value = 69.29994029850732
$ws.Cells.Item(47, 2).Value = 53.29999999999999
$ws.Cells.Item(47, 3).Value = 82.16400000000007
$ws.Cells.Item(48, 1).Value = 74.05586069651724
$ws.Cells.Item(48, 2).Value = 60.16000000000004
$ws.Cells.Item(48, 3).Value = 84.19999999999996
$ws.Cells.Item(49, 1).Value = 42.15331343283577
$ws.Cells.Item(49, 2).Value = 27.01600000000001
$ws.Cells.Item(49, 3).Value = 58.42
$ws.Cells.Item(50, 1).Value = 51.97609950248741
$ws.Cells.Item(50, 2).Value = 35.82000000000001
$ws.Cells.Item(50, 3).Value = 67.18799999999997
$ws.Cells.Item(51, 1).Value = 26.21671641791042
$ws.Cells.Item(51, 2).Value = 15.16
$ws.Cells.Item(51, 3).Value = 39.61999999999998
$ws.Cells.Item(52, 1).Value = 76.14712437810927
$ws.Cells.Item(52, 2).Value = 62.15600000000004
$ws.Cells.Item(52, 3).Value = 86.33600000000004
$ws.Cells.Item(53, 1).Value = 75.08348258706467
$ws.Cells.Item(53, 2).Value = 62.25200000000002
$ws.Cells.Item(53, 3).Value = 85.96799999999996
$ws.Cells.Item(54, 1).Value = 33.85383084577103
$ws.Cells.Item(54, 2).Value = 20.34799999999998
$ws.Cells.Item(54, 3).Value = 50.21600000000002
$ws.Cells.Item(55, 1).Value = 48.55203980099496
$ws.Cells.Item(55, 2).Value = 30.87199999999999
$ws.Cells.Item(55, 3).Value = 65.17200000000001
$ws.Cells.Item(56, 1).Value = 77.36630845771136
$ws.Cells.Item(56, 2).Value = 65.3
$ws.Cells.Item(56, 3).Value = 87.08399999999995
$ws.Cells.Item(57, 1).Value = 73.58754228855707
$ws.Cells.Item(57, 2).Value = 60.65599999999996
$ws.Cells.Item(57, 3).Value = 85.18400000000001
$ws.Cells.Item(58, 1).Value = 51.89769154228846
$ws.Cells.Item(58, 2).Value = 33.336
$ws.Cells.Item(58, 3).Value = 67.508
$ws.Cells.Item(59, 1).Value = 77.63317412935318
$ws.Cells.Item(59, 2).Value = 66.40399999999998
$ws.Cells.Item(59, 3).Value = 87.02399999999997
$ws.Cells.Item(60, 1).Value = 37.0556019900496
$ws.Cells.Item(60, 2).Value = 23.396
$ws.Cells.Item(60, 3).Value = 53.4
$ws.Cells.Item(61, 1).Value = 77.39844776119389
$ws.Cells.Item(61, 2).Value = 65.03999999999998
$ws.Cells.Item(61, 3).Value = 87.7
$ws.Cells.Item(62, 1).Value = 78.68583084577108
$ws.Cells.Item(62, 2).Value = 66.75599999999997
$ws.Cells.Item(62, 3).Value = 87.63599999999995
$ws.Cells.Item(63, 1).Value = 78.58847761194021
$ws.Cells.Item(63, 2).Value = 68.29599999999995
$ws.Cells.Item(63, 3).Value = 88.05199999999999
$ws.Cells.Item(64, 1).Value = 37.44392039800983
$ws.Cells.Item(64, 2).Value = 24.228
$ws.Cells.Item(64, 3).Value = 53.71199999999999
$ws.Cells.Item(65, 1).Value = 73.46811940298497
$ws.Cells.Item(65, 2).Value = 59.86000000000001
$ws.Cells.Item(65, 3).Value = 84.892
$ws.Cells.Item(66, 1).Value = 78.11880597014925
$ws.Cells.Item(66, 2).Value = 66.508
$ws.Cells.Item(66, 3).Value = 86.93599999999996
$ws.Cells.Item(67, 1).Value = 78.49454726368151
$ws.Cells.Item(67, 2).Value = 66.84799999999997
$ws.Cells.Item(67, 3).Value = 88.24000000000004
$ws.Cells.Item(68, 1).Value = 47.94997014925362
$ws.Cells.Item(68, 2).Value = 30.48400000000002
$ws.Cells.Item(68, 3).Value = 65.304
$ws.Cells.Item(69, 1).Value = 69.70437810945268
$ws.Cells.Item(69, 2).Value = 54.29600000000001
$ws.Cells.Item(69, 3).Value = 83.196
$ws.Cells.Item(70, 1).Value = 44.64565174129342
$ws.Cells.Item(70, 2).Value = 29.71600000000002
$ws.Cells.Item(70, 3).Value = 61.98800000000001
$ws.Cells.Item(71, 1).Value = 65.19494527363173
$ws.Cells.Item(71, 2).Value = 47.89999999999999
$ws.Cells.Item(71, 3).Value = 78.44399999999996
$ws.Cells.Item(72, 1).Value = 70.46714427860688
$ws.Cells.Item(72, 2).Value = 56.07199999999995
$ws.Cells.Item(72, 3).Value = 82.56800000000003
